$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: header text "Geopotential height(hpm)" -> "GH", centered (no vertical-top) alignment
$ws.Range("B1").Value = "GH"
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter

# B8: -999 -> text "nan", right aligned
$ws.Range("B8").Value = "nan"
$ws.Range("B8").HorizontalAlignment = -4152   # xlRight

# B29: 556 -> text "nan", right aligned
$ws.Range("B29").Value = "nan"
$ws.Range("B29").HorizontalAlignment = -4152  # xlRight

# Update active selection from D5 to D9
$ws.Range("D9").Select()
